$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple D/E value updates (rows unaffected by reordering) ---
$ws.Cells.Item(2, 4).Value = "44.016.18"
$ws.Cells.Item(2, 5).Value = "  +0.86%  "

$ws.Cells.Item(3, 4).Value = "2.331.03"
$ws.Cells.Item(3, 5).Value = "  +4.23%  "

$ws.Cells.Item(4, 5).Value = "  +0.10%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "97.57"
$ws.Cells.Item(5, 5).Value = "  +4.19%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "271.95"
$ws.Cells.Item(6, 5).Value = "  +0.96%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.627"
$ws.Cells.Item(7, 5).Value = "  +0.51%  "

$ws.Cells.Item(8, 5).Value = "  +0.08%  "

$ws.Cells.Item(9, 5).Value = "  +0.89%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "46.29"
$ws.Cells.Item(10, 5).Value = "  -1.04%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0956"
$ws.Cells.Item(11, 5).Value = "  +3.43%  "

$ws.Cells.Item(12, 5).Value = "  -2.21%  "

$ws.Cells.Item(13, 5).Value = "  +0.05%  "

$ws.Cells.Item(14, 4).Value = "2.682.57"
$ws.Cells.Item(14, 5).Value = "  +4.40%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "15.65"
$ws.Cells.Item(15, 5).Value = "  +3.35%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.873"
$ws.Cells.Item(16, 5).Value = "  +8.96%  "

$ws.Cells.Item(17, 4).Value = "2.331.50"
$ws.Cells.Item(17, 5).Value = "  +3.96%  "

$ws.Cells.Item(18, 4).Value = "43.908.12"
$ws.Cells.Item(18, 5).Value = "  +0.82%  "

$ws.Cells.Item(19, 5).Value = "  +5.44%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "6.45"
$ws.Cells.Item(20, 5).Value = "  +7.21%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "72.94"
$ws.Cells.Item(21, 5).Value = "  +3.48%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "240.17"
$ws.Cells.Item(22, 5).Value = "  +2.86%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "2.29"
$ws.Cells.Item(23, 5).Value = "  -1.93%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "9.49"
$ws.Cells.Item(24, 5).Value = "  +5.08%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "0.999"
$ws.Cells.Item(25, 5).Value = "  -0.08%  "

$ws.Cells.Item(28, 5).Value = "  -1.88%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "2.28"
$ws.Cells.Item(29, 5).Value = "  +0.47%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "38.38"
$ws.Cells.Item(30, 5).Value = "  -5.08%  "

$ws.Cells.Item(31, 5).Value = "  +7.79%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "173.82"
$ws.Cells.Item(32, 5).Value = "  +0.43%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.0908"
$ws.Cells.Item(33, 5).Value = "  -2.15%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "5.52"
$ws.Cells.Item(34, 5).Value = "  +0.71%  "

$ws.Cells.Item(35, 5).Value = "  +2.58%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.0363"
$ws.Cells.Item(36, 5).Value = "  +3.56%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.109"
$ws.Cells.Item(37, 5).Value = "  -1.77%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "4.45"
$ws.Cells.Item(38, 5).Value = "  +2.72%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "3.39"
$ws.Cells.Item(39, 5).Value = "  -6.17%  "

$ws.Cells.Item(42, 5).Value = "  +18.44%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "12.30"
$ws.Cells.Item(43, 5).Value = "  -2.38%  "

$ws.Cells.Item(44, 5).Value = "  +9.73%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "62.50"
$ws.Cells.Item(45, 5).Value = "  -1.27%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "5.42"
$ws.Cells.Item(46, 5).Value = "  +1.48%  "

$ws.Cells.Item(47, 5).Value = "  +4.67%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "100.50"
$ws.Cells.Item(48, 5).Value = "  -0.34%  "

$ws.Cells.Item(49, 5).Value = "  +1.49%  "

$ws.Cells.Item(50, 4).Value = "2.559.68"
$ws.Cells.Item(50, 5).Value = "  +4.25%  "

$ws.Cells.Item(51, 5).Value = "  +15.90%  "

# --- Row 26/27 swap: Cosmos <-> PancakeSwap (values refreshed) ---
$ws.Cells.Item(26, 2).Value = "PancakeSwap"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "2.54"
$ws.Cells.Item(26, 5).Value = "  +1.30%  "

$ws.Cells.Item(27, 2).Value = "Cosmos"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "11.43"
$ws.Cells.Item(27, 5).Value = "  +1.17%  "

# --- Row 40/41 swap: LidoDAOToken <-> Algorand (values refreshed) ---
$ws.Cells.Item(40, 2).Value = "Algorand"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.240"
$ws.Cells.Item(40, 5).Value = "  +9.47%  "

$ws.Cells.Item(41, 2).Value = "LidoDAOToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "2.37"
$ws.Cells.Item(41, 5).Value = "  +8.85%  "
